$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- O13: change cell format (border + centered, no longer vertically centered) ---
# Matches the pre-existing "style 1" cell format (thin border + horizontal-center only)
$ws.Range("O13").HorizontalAlignment = -4108
$ws.Range("O13").Borders.LineStyle = 1

# --- New data rows 43-52 in columns X:Z (continuation of the viscosity table) ---
$newRows = @(
    @{ Row = 43; X = 10;  Y = 1446.55; Z = 3.14 },
    @{ Row = 44; X = 20;  Y = 835.66;  Z = 2.59 },
    @{ Row = 45; X = 30;  Y = 704.13;  Z = 2.93 },
    @{ Row = 46; X = 40;  Y = 405.08;  Z = 2.76 },
    @{ Row = 47; X = 50;  Y = 316.49;  Z = 2.5 },
    @{ Row = 48; X = 60;  Y = 343.96;  Z = 2.78 },
    @{ Row = 49; X = 70;  Y = 1319.71; Z = 2.81 },
    @{ Row = 50; X = 80;  Y = $null;   Z = $null },
    @{ Row = 51; X = 90;  Y = $null;   Z = $null },
    @{ Row = 52; X = 100; Y = $null;   Z = $null }
)

foreach ($r in $newRows) {
    $xCell = $ws.Cells.Item($r.Row, 24)
    $ws.Range("X39").Copy()
    $xCell.PasteSpecial(-4122)
    $xCell.Value = $r.X

    if ($null -ne $r.Y) {
        $ws.Cells.Item($r.Row, 25).Value = $r.Y
    }
    if ($null -ne $r.Z) {
        $ws.Cells.Item($r.Row, 26).Value = $r.Z
    }
}

$excel.CutCopyMode = $false

# --- Sheet view: zoom + scroll position + active selection ---
$excel.ActiveWindow.Zoom = 102
$ws.Range("Y48").Select()
